$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D; old D:K shift to E:L
$ws.Range("D1").EntireColumn.Insert()

# 2) Copy column E formatting (number format/style) into the new column D
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Match column D width to its neighbor
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# 4) Write the full D:K data grid (new column D plus corrected historical
#    restatements that landed in E/F after the shift)
$rowData = @{
    7 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
    8 = @(1082800, 1236900, 2302100, 3352300, 3232500, 2538100, 2200700, 2695800)
    9 = @(667000, 661400, 923200, 1302800, 1566900, 1230800, 1047700, 1476500)
    10 = @(415800, 575500, 1378900, 2049400, 1665600, 1307400, 1153000, 1219300)
    11 = @($null, $null, $null, $null, $null, $null, $null, $null)
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    13 = @(0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(803900, 121600, 1440900, 418300, 745400, -9300, -12900, -21200)
    15 = @(486500, 548000, 611100, 634300, 627500, 511500, 440300, 658600)
    16 = @($null, $null, $null, $null, $null, $null, $null, $null)
    17 = @(2030700, 1402600, 3044400, 2432300, 3046600, 1797600, 1575100, 2205300)
    18 = @(-947900, -165700, -742400, 920000, 185900, 740500, 625600, 490500)
    19 = @($null, $null, $null, $null, $null, $null, $null, $null)
    20 = @(8300, 7900, -1700, 36300, -1300, 4200, 3600, 1500)
    21 = @(-453000, 390200, -133000, 1590600, 1048200, 1624100, 1387700, 1150600)
    22 = @(297600, 292000, 222900, 213900, 155200, 106300, 85800, 55700)
    23 = @(-1237200, -449800, -967000, 742400, 29500, 638400, 543400, 436300)
    24 = @(-106600, 151600, -109200, 159200, 106700, 92100, 95200, 72600)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(-1130500, -601400, -857900, 583200, -77200, 546300, 448200, 363600)
    27 = @(-885100, -624000, -929600, 500100, -152000, 469300, 409100, 367200)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0)
    29 = @(0, 107500, 0, 0, 160500, 304100, 108000, "NA")
    30 = @(0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(-8300, -7900, 1700, -36300, 1300, -4200, -3600, -1500)
    33 = @(-885100, -516500, -929600, 500100, 8500, 773400, 517000, 367200)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(-885100, -516500, -929600, 500100, 8500, 773400, 517000, 367200)
    38 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
    39 = @($null, $null, $null, $null, $null, $null, $null, $null)
    40 = @($null, $null, $null, $null, $null, $null, $null, $null)
    41 = @(375200, 662800, 725700, 512200, 68500, 114500, 282100, 239200)
    42 = @(0, 0, 0, 0, 0, 0, 0, 0)
    43 = @(221200, 310000, 374600, 554500, 676600, 1089300, 856100, 662400)
    44 = @(0, 0, 0, 0, 0, 0, 0, 0)
    45 = @(62600, 66100, 92300, 173900, 183500, 187100, 167100, 158000)
    46 = @(659100, 1039000, 1192600, 1240600, 928600, 1390900, 1305300, 1059600)
    47 = @(0, 0, 0, 0, 0, 0, 0, 0)
    48 = @(8480700, 9489200, 10061900, 11483600, 12112500, 14558100, 13026000, 12130300)
    49 = @(0, 0, 0, 0, 0, 0, 0, 0)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(125100, 266400, 185600, 141400, 245800, 268900, 276500, 305200)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(9264900, 10794700, 11440100, 12865600, 13286800, 16218000, 14607800, 13495200)
    55 = @($null, $null, $null, $null, $null, $null, $null, $null)
    56 = @($null, $null, $null, $null, $null, $null, $null, $null)
    57 = @(125600, 84000, 108200, 223200, 265400, 347200, 350100, 436000)
    58 = @(0, 249800, 299900, 299900, "NA", "NA", "NA", 0)
    59 = @(239900, 259100, 225200, 340400, 403300, 704700, 561300, 515100)
    60 = @(365500, 593000, 633300, 863600, 668700, 1051900, 911400, 827200)
    61 = @(3877400, 3795900, 4040200, 4162600, 4869000, 5556300, 4634400, 4072000)
    62 = @(367500, 455100, 299100, 417200, 462100, 559800, 573700, 548300)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(5011800, 5518500, 5681400, 6166400, 6722100, 7895400, 6884600, 6088600)
    67 = @($null, $null, $null, $null, $null, $null, $null, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(3608400, 4637700, 5154200, 6131500, 5936000, 7591900, 7066000, 6676400)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(4253200, 5276200, 5758700, 6699200, 6564700, 8322600, 7723200, 7406500)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0)
    80 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
    81 = @(-885100, -516500, -929600, 500100, 8500, 773400, 517000, 367200)
    82 = @($null, $null, $null, $null, $null, $null, $null, $null)
    83 = @(486500, 548000, 611100, 634300, 863500, 879400, 758600, 658600)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(171900, 416700, 1142700, 1764900, 1778200, 1702300, 1381700, 740200)
    90 = @($null, $null, $null, $null, $null, $null, $null, $null)
    91 = @(-194800, -120700, -711400, -422500, -2072900, -2487500, -1669800, -2621200)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-189400, -118300, -686600, -432500, -2109300, -2485100, -1790900, -2521500)
    95 = @($null, $null, $null, $null, $null, $null, $null, $null)
    96 = @(0, 0, -47500, -315500, -386600, -194900, -138300, -150500)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(-269400, -361200, -242700, -888600, 285100, 615200, 452100, 1682600)
    101 = @(0, 0, 0, 0, 0, 0, 0, 0)
    102 = @(-286900, -62900, 213500, 443700, -45900, -167600, 42900, -98700)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 4 + $i   # D=4 .. K=11
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item([int]$r, $col).Value = $v
        }
    }
}

Write-Host "done"